$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (F value, G value) updates for F221:G270
$data = @(
    @(221, 0, 0),
    @(222, 525, 40),
    @(223, 1152, 84),
    @(224, 929, 31),
    @(225, 372, 36),
    @(226, 320, 27),
    @(227, 0, 0),
    @(228, 92, 0),
    @(229, 792, 69),
    @(230, 614, 38),
    @(231, 954, 48),
    @(232, 2138, 141),
    @(233, 2161, 119),
    @(234, 1194, 65),
    @(235, 729, 40),
    @(236, 2161, 231),
    @(237, 3036, 234),
    @(238, 2415, 225),
    @(239, 6185, 239),
    @(240, 41366, 492),
    @(241, 89674, 975),
    @(242, 30176, 376),
    @(243, 26701, 1331),
    @(244, 5283, 98),
    @(245, 3606, 82),
    @(246, 2115, 140),
    @(247, 34523, 327),
    @(248, 45188, 439),
    @(249, 12172, 143),
    @(250, 17656, 1154),
    @(251, 6492, 140),
    @(252, 4620, 101),
    @(253, 5636, 173),
    @(254, 6145, 158),
    @(255, 5452, 99),
    @(256, 1245, 33),
    @(257, 5642, 286),
    @(258, 3919, 224),
    @(259, 6546, 442),
    @(260, 12193, 770),
    @(261, 18110, 611),
    @(262, 9285, 322),
    @(263, 1701, 57),
    @(264, 43207, 885),
    @(265, 18446, 894),
    @(266, 14726, 733),
    @(267, 15446, 867),
    @(268, 17529, 787),
    @(269, 9705, 433),
    @(270, 2780, 180)
)

foreach ($item in $data) {
    $row = $item[0]
    $fValue = $item[1]
    $gValue = $item[2]
    $ws.Cells.Item($row, 6).Value = $fValue
    $ws.Cells.Item($row, 7).Value = $gValue
}

Write-Output "Updated F221:G270 with AgTests/AgPosit values."
